# Apply corrected error-estimation / projected-years results to the
# SoIB_summaries workbook (Madhya Pradesh).

$wb = $excel.ActiveWorkbook

# --- "Trends Status" sheet ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("C2").Value = 3
$ws1.Range("E2").Value = 5.4
$ws1.Range("C3").Value = 18
$ws1.Range("E3").Value = 32.1
$ws1.Range("C4").Value = 29
$ws1.Range("E4").Value = 51.8
$ws1.Range("E5").Value = 3.6
$ws1.Range("E6").Value = 7.1
$ws1.Range("C7").Value = 48

# --- "Species qualification" sheet ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("C4").Value = 56

# --- "Interannual update - High Pri" sheet ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")
$ws5.Range("B2").Value = 74
$ws5.Range("C2").Value = 71.8
$ws5.Range("D2").Value = 74
$ws5.Range("E2").Value = 79.59999999999999
$ws5.Range("B3").Value = 29
$ws5.Range("C3").Value = 28.2
$ws5.Range("D3").Value = 19
$ws5.Range("E3").Value = 20.4
